$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Severity-Mortality")
$ws.Activate()
